$wb = $excel.ActiveWorkbook

$urlA = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5a80c0ec42a04012333b4102b656fc428d33ee91/e2e/831a32a9-0652-4582-a26e-0776f09d60aa.md"
$urlB = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5a80c0ec42a04012333b4102b656fc428d33ee91/e2e/e4723d53-a48f-4b2c-8181-85f1f3fc95db.md"
$displayA = "831a32a9-0652-4582-a26e-0776f09d60aa.md"
$displayB = "e4723d53-a48f-4b2c-8181-85f1f3fc95db.md"

# The handback run finished for both localized files - flip the Overview/per-locale
# "Status" column from "Ready for handoff" to the synced handback message.
# (Status is column C on the per-locale sheets and shares a string with
# the Overview sheet's zh-cn/de-de columns, so replace it everywhere it's used
# instead of just overwriting individual cells - that keeps every occurrence
# pointed at a single shared string, same as the source text did.)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

$wsOverview.Cells.Replace("Ready for handoff", "Handed back: in sync with en-US")
$wsZh.Cells.Replace("Ready for handoff", "Handed back: in sync with en-US")
$wsDe.Cells.Replace("Ready for handoff", "Handed back: in sync with en-US")

# Fill in the generated handback report for zh-cn: target file, handback file
# and handback datetime for both localized documents, with the target file
# name re-hyperlinked to the source doc (matching column A's link style).
$wsZh.Range("I2").Value = $displayA
$wsZh.Range("J2").Value = "831a32a9-0652-4582-a26e-0776f09d60aa.7e587a0d3f96732684d9117407d33c41bb5c9eba.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-19 00:29:32"

$wsZh.Range("I3").Value = $displayB
$wsZh.Range("J3").Value = "e4723d53-a48f-4b2c-8181-85f1f3fc95db.cd4c01b1c96703a688d36c74cbfab560823ffe3d.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-19 00:29:32"

# Re-create the two existing doc-name hyperlinks so the new Latest-Target-File
# hyperlinks slot in right after their row's first link (keeps rIds in row order).
$wsZh.Range("A2").Hyperlinks.Delete()
$wsZh.Range("A3").Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $urlA, "", "", $displayA)
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $urlA, "", "", $displayA)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $urlB, "", "", $displayB)
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $urlB, "", "", $displayB)

# Same handback report fill-in for de-de.
$wsDe.Range("C2:C3").Value = "Handed back: in sync with en-US"

$wsDe.Range("I2").Value = $displayA
$wsDe.Range("J2").Value = "831a32a9-0652-4582-a26e-0776f09d60aa.7e587a0d3f96732684d9117407d33c41bb5c9eba.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-19 00:29:39"

$wsDe.Range("I3").Value = $displayB
$wsDe.Range("J3").Value = "e4723d53-a48f-4b2c-8181-85f1f3fc95db.cd4c01b1c96703a688d36c74cbfab560823ffe3d.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-19 00:29:39"

$wsDe.Range("A2").Hyperlinks.Delete()
$wsDe.Range("A3").Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $urlA, "", "", $displayA)
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $urlA, "", "", $displayA)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $urlB, "", "", $displayB)
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $urlB, "", "", $displayB)

# Widen the Status column (now holding the longer handback message) and the
# Latest Target File / Latest Handback File columns (now holding filenames)
# to match Excel's autofit after the content grew.
$wsZh.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsZh.Columns.Item(9).ColumnWidth = 40
$wsZh.Columns.Item(10).ColumnWidth = 40

$wsDe.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsDe.Columns.Item(9).ColumnWidth = 40
$wsDe.Columns.Item(10).ColumnWidth = 40

$wsOverview.Columns.Item(5).ColumnWidth = 29.9777047293527
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777047293527

Write-Output "handback report generated"
